$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 12) down to the
# two new rows (13 and 14) before/after filling in values, so the new
# cells pick up the same number format (date) / borders as the rest of
# the table.
$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B12:M12").Copy()
$ws.Range("B13:M14").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Row 13 data
$ws.Range("A13").Value = 45762
$ws.Range("B13").Value = 36.9
$ws.Range("C13").Value = 42.9
$ws.Range("D13").Value = 38.7
$ws.Range("E13").Value = 40.3
$ws.Range("F13").Value = 38.9
$ws.Range("G13").Value = 32.9
$ws.Range("H13").Value = 33
$ws.Range("I13").Value = 38.1
$ws.Range("J13").Value = 28.5
$ws.Range("K13").Value = 29.6
$ws.Range("L13").Value = 25.6
$ws.Range("M13").Value = 30.8

# Row 14 data
$ws.Range("A14").Value = 45763
$ws.Range("B14").Value = 35.7
$ws.Range("C14").Value = 40.9
$ws.Range("D14").Value = 38
$ws.Range("E14").Value = 39.7
$ws.Range("F14").Value = 37.6
$ws.Range("G14").Value = 32.1
$ws.Range("H14").Value = 31.5
$ws.Range("I14").Value = 37.4
$ws.Range("J14").Value = 28.1
$ws.Range("K14").Value = 28.7
$ws.Range("L14").Value = 25.9
$ws.Range("M14").Value = 30.1

# Update selection to match the saved view state
$ws.Range("E16").Select()
